$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Insert a new weekly data row above row 132; this shifts the existing
# rows 132:155 down to 133:156 (dimension grows from R155 to R156).
$ws.Rows.Item(132).Insert()

# Populate the newly inserted row 132 with the latest weekly record.
$ws.Range("A132").Value = 7
$ws.Range("B132").Value = "Terminal Hortofrutícola Agro Chillán"
$ws.Range("C132").Value = "Ñuble"
$ws.Range("D132").Value = 45124
$ws.Range("E132").Value = 16
$ws.Range("F132").Value = 100112031
$ws.Range("G132").Value = "Poroto verde"
$ws.Range("H132").Value = "Magnum"
$ws.Range("I132").Value = "Primera"
$ws.Range("J132").Value = 30
$ws.Range("K132").Value = 25000
$ws.Range("L132").Value = 25000
$ws.Range("M132").Value = 25000
$ws.Range("N132").Value = "$/malla 25 kilos"
$ws.Range("O132").Value = "Perú"
$ws.Range("P132").Value = 1000
$ws.Range("Q132").Value = 25
$ws.Range("R132").Value = "Hortaliza"
